$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 4 for columns D, M, N, O, P, Q, S, T
# (row 2 becomes what row 4 was, and vice versa)

$ws.Range("D2").Value = 44875
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("S2").Value = 1600
$ws.Range("T2").Value = 10

$ws.Range("D4").Value = 44855
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/bandeja 5 kilos"
$ws.Range("S4").Value = 3000
$ws.Range("T4").Value = 5
